$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '37.507.40'
$ws.Range("E2").Value = '  +0.92%  '

$ws.Range("D3").Value = '2.075.89'
$ws.Range("E3").Value = '  +3.71%  '

$ws.Range("E4").Value = '  +0.08%  '

$ws.Range("D5").Value = '235.84'
$ws.Range("E5").Value = '  -3.67%  '

$ws.Range("D6").Value = '0.617'
$ws.Range("E6").Value = '  +1.62%  '

$ws.Range("D7").Value = '57.78'
$ws.Range("E7").Value = '  +3.73%  '

$ws.Range("E8").Value = '  +0.03%  '

$ws.Range("D9").Value = '0.383'
$ws.Range("E9").Value = '  +1.98%  '

$ws.Range("D10").Value = '58.14'
$ws.Range("E10").Value = '  +0.07%  '

$ws.Range("D11").Value = '0.0762'
$ws.Range("E11").Value = '  +0.24%  '

$ws.Range("D12").Value = '0.101'
$ws.Range("E12").Value = '  +3.35%  '

$ws.Range("D13").Value = '2.382.72'
$ws.Range("E13").Value = '  +4.18%  '

$ws.Range("D14").Value = '14.55'
$ws.Range("E14").Value = '  +2.28%  '

$ws.Range("D15").Value = '21.11'
$ws.Range("E15").Value = '  +0.09%  '

$ws.Range("D16").Value = '0.777'
$ws.Range("E16").Value = '  +1.58%  '

$ws.Range("D17").Value = '5.25'
$ws.Range("E17").Value = '  +3.57%  '

$ws.Range("D18").Value = '2.051.46'
$ws.Range("E18").Value = '  +2.47%  '

$ws.Range("D19").Value = '37.633.67'
$ws.Range("E19").Value = '  +1.55%  '

$ws.Range("D20").Value = '5.98'
$ws.Range("E20").Value = '  +18.92%  '

$ws.Range("D21").Value = '68.46'
$ws.Range("E21").Value = '  +0.03%  '

$ws.Range("D22").Value = '0.0₃0812'
$ws.Range("E22").Value = '  -0.31%  '

$ws.Range("D23").Value = '223.98'
$ws.Range("E23").Value = '  -2.42%  '

$ws.Range("E24").Value = '  -0.07%  '

$ws.Range("E25").Value = '  +3.30%  '

$ws.Range("D26").Value = '2.44'
$ws.Range("E26").Value = '  -0.93%  '

$ws.Range("D27").Value = '162.86'
$ws.Range("E27").Value = '  +0.26%  '

$ws.Range("D28").Value = '8.88'
$ws.Range("E28").Value = '  +1.89%  '

$ws.Range("D29").Value = '0.132'
$ws.Range("E29").Value = '  +5.18%  '

$ws.Range("D30").Value = '19.34'
$ws.Range("E30").Value = '  +0.56%  '

$ws.Range("E31").Value = '  +5.27%  '

$ws.Range("E32").Value = '  +0.87%  '

$ws.Range("D33").Value = '4.47'
$ws.Range("E33").Value = '  +0.23%  '

$ws.Range("D34").Value = '0.0625'
$ws.Range("E34").Value = '  +1.70%  '

$ws.Range("D35").Value = '2.57'
$ws.Range("E35").Value = '  +9.62%  '

$ws.Range("D36").Value = '4.39'
$ws.Range("E36").Value = '  +2.85%  '

$ws.Range("E37").Value = '  +0.04%  '

$ws.Range("B38").Value = 'RenderToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D38").Value = '3.36'
$ws.Range("E38").Value = '  -0.17%  '

$ws.Range("B39").Value = 'THORChain'
$ws.Range("C39").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D39").Value = '5.94'
$ws.Range("E39").Value = '  +13.14%  '

$ws.Range("D40").Value = '1.77'
$ws.Range("E40").Value = '  -0.91%  '

$ws.Range("E41").Value = '  -4.77%  '

$ws.Range("D42").Value = '0.0966'
$ws.Range("E42").Value = '  +9.03%  '

$ws.Range("D43").Value = '1.475.79'
$ws.Range("E43").Value = '  +2.39%  '

$ws.Range("D44").Value = '4.33'
$ws.Range("E44").Value = '  +22.19%  '

$ws.Range("D45").Value = '95.17'
$ws.Range("E45").Value = '  +7.10%  '

$ws.Range("D46").Value = '16.51'
$ws.Range("E46").Value = '  +7.09%  '

$ws.Range("E47").Value = '  +2.65%  '

$ws.Range("E48").Value = '  +0.21%  '

$ws.Range("E49").Value = '  +1.97%  '

$ws.Range("E50").Value = '  +9.21%  '

$ws.Range("E51").Value = '  +1.08%  '
